# Weekly update: append new "Sandia" price rows for Vega Monumental Concepción.
#
# The sheet is a flat data table (header in row 1, data from row 2 on).
# This week's update inserts three new observation rows:
#   - two new rows at the very top of the existing "Sandia" block (new rows 55-56)
#   - one new row further down, right after the old row that is now row 66
#     (new row 67)
# All previously existing rows simply shift down to make room; their contents
# are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the first two new rows at row 55 -------------------------------
$ws.Rows(55).Insert()
$ws.Rows(55).Insert()

$ws.Range("A55").Value = 11
$ws.Range("B55").Value = "Vega Monumental Concepción"
$ws.Range("C55").Value = "Bíobío"
$ws.Range("D55").Value = 44587
$ws.Range("E55").Value = 8
$ws.Range("F55").Value = 100112028
$ws.Range("G55").Value = "Sandia"
$ws.Range("H55").Value = "Sin especificar"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 1800
$ws.Range("K55").Value = 2000
$ws.Range("L55").Value = 2300
$ws.Range("M55").Value = 2133
$ws.Range("N55").Value = '$/unidad'
$ws.Range("O55").Value = "Región de O'Higgins"
$ws.Range("P55").Value = 2133
$ws.Range("Q55").Value = 1
$ws.Range("R55").Value = "Hortaliza"

$ws.Range("A56").Value = 11
$ws.Range("B56").Value = "Vega Monumental Concepción"
$ws.Range("C56").Value = "Bíobío"
$ws.Range("D56").Value = 44587
$ws.Range("E56").Value = 8
$ws.Range("F56").Value = 100112028
$ws.Range("G56").Value = "Sandia"
$ws.Range("H56").Value = "Sin especificar"
$ws.Range("I56").Value = "Segunda"
$ws.Range("J56").Value = 2000
$ws.Range("K56").Value = 1400
$ws.Range("L56").Value = 1500
$ws.Range("M56").Value = 1450
$ws.Range("N56").Value = '$/unidad'
$ws.Range("O56").Value = "Región de O'Higgins"
$ws.Range("P56").Value = 1450
$ws.Range("Q56").Value = 1
$ws.Range("R56").Value = "Hortaliza"

# --- Insert the third new row at row 67 (after the data already shifted) --
$ws.Rows(67).Insert()

$ws.Range("A67").Value = 11
$ws.Range("B67").Value = "Vega Monumental Concepción"
$ws.Range("C67").Value = "Bíobío"
$ws.Range("D67").Value = 44511
$ws.Range("E67").Value = 8
$ws.Range("F67").Value = 100112028
$ws.Range("G67").Value = "Sandia"
$ws.Range("H67").Value = "Sin especificar"
$ws.Range("I67").Value = "Segunda"
$ws.Range("J67").Value = 300
$ws.Range("K67").Value = 700
$ws.Range("L67").Value = 700
$ws.Range("M67").Value = 700
$ws.Range("N67").Value = '$/kilo (volumen en unidades)'
$ws.Range("O67").Value = "Perú"
$ws.Range("P67").Value = 700
$ws.Range("Q67").Value = 1
$ws.Range("R67").Value = "Hortaliza"
